# Refresh the cryptocurrency table (columns B:E, rows 2-51) with the
# newer Coinranking snapshot values.
#
# Commit message: "Updated cryptos list on Fri Nov 10 00:56:58 UTC 2023
# with GitHub Actions"
#
# A handful of rows (21/22, 32/33/34, 35/37, 42/43, 51) also changed
# which coin occupies that rank, so their Coin name + Link are updated
# along with Price/Volume(1h).
#
# The Price/Volume(1h) columns hold plain text such as "36.533.46" or
# "  +1.90%  " (note the deliberate padding + thousands-as-dots
# formatting). Assigning those strings straight to .Value would let
# Excel's smart-parsing reinterpret many of them as numbers/percentages
# and silently mangle the text (e.g. "253.00" -> 253, "0.0978" ->
# 9.7800000000000005E-2). Forcing the cells to Text format first, then
# restoring the Normal style afterwards, keeps the original text intact
# without leaving a custom number format applied to the cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$priceVolRange = $ws.Range("D2:E51")
$priceVolRange.NumberFormat = "@"

$ws.Cells.Item(2, 4).Value = '36.533.46'
$ws.Cells.Item(2, 5).Value = '  +1.90%  '
$ws.Cells.Item(3, 4).Value = '2.100.90'
$ws.Cells.Item(3, 5).Value = '  +10.77%  '
$ws.Cells.Item(4, 5).Value = '  +0.08%  '
$ws.Cells.Item(5, 4).Value = '253.00'
$ws.Cells.Item(5, 5).Value = '  +2.29%  '
$ws.Cells.Item(6, 4).Value = '0.664'
$ws.Cells.Item(6, 5).Value = '  -3.63%  '
$ws.Cells.Item(7, 5).Value = '  +0.03%  '
$ws.Cells.Item(8, 4).Value = '45.42'
$ws.Cells.Item(8, 5).Value = '  +5.82%  '
$ws.Cells.Item(9, 4).Value = '61.78'
$ws.Cells.Item(9, 5).Value = '  +9.30%  '
$ws.Cells.Item(10, 5).Value = '  +1.93%  '
$ws.Cells.Item(11, 4).Value = '0.0737'
$ws.Cells.Item(11, 5).Value = '  -2.92%  '
$ws.Cells.Item(12, 4).Value = '0.0996'
$ws.Cells.Item(12, 5).Value = '  +0.80%  '
$ws.Cells.Item(13, 4).Value = '14.56'
$ws.Cells.Item(13, 5).Value = '  -1.53%  '
$ws.Cells.Item(14, 4).Value = '2.409.95'
$ws.Cells.Item(14, 5).Value = '  +10.86%  '
$ws.Cells.Item(15, 4).Value = '0.843'
$ws.Cells.Item(15, 5).Value = '  +6.20%  '
$ws.Cells.Item(16, 4).Value = '2.117.21'
$ws.Cells.Item(16, 5).Value = '  +11.59%  '
$ws.Cells.Item(17, 5).Value = '  +0.57%  '
$ws.Cells.Item(18, 4).Value = '36.674.34'
$ws.Cells.Item(18, 5).Value = '  +2.26%  '
$ws.Cells.Item(19, 4).Value = '73.82'
$ws.Cells.Item(19, 5).Value = '  +0.57%  '
$ws.Cells.Item(20, 4).Value = '0.0₃0827'
$ws.Cells.Item(20, 5).Value = '  -0.68%  '
$ws.Cells.Item(21, 2).Value = 'Avalanche'
$ws.Cells.Item(21, 3).Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Cells.Item(21, 4).Value = '13.06'
$ws.Cells.Item(21, 5).Value = '  +0.10%  '
$ws.Cells.Item(22, 2).Value = 'BitcoinCash'
$ws.Cells.Item(22, 3).Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Cells.Item(22, 4).Value = '240.07'
$ws.Cells.Item(22, 5).Value = '  -2.89%  '
$ws.Cells.Item(23, 4).Value = '5.13'
$ws.Cells.Item(23, 5).Value = '  -1.54%  '
$ws.Cells.Item(24, 5).Value = '  +0.09%  '
$ws.Cells.Item(25, 4).Value = '2.46'
$ws.Cells.Item(25, 5).Value = '  -7.84%  '
$ws.Cells.Item(26, 4).Value = '169.83'
$ws.Cells.Item(26, 5).Value = '  +1.85%  '
$ws.Cells.Item(27, 4).Value = '21.37'
$ws.Cells.Item(27, 5).Value = '  +15.88%  '
$ws.Cells.Item(28, 4).Value = '9.20'
$ws.Cells.Item(28, 5).Value = '  +5.40%  '
$ws.Cells.Item(29, 4).Value = '2.00'
$ws.Cells.Item(29, 5).Value = '  -9.51%  '
$ws.Cells.Item(30, 4).Value = '0.123'
$ws.Cells.Item(30, 5).Value = '  -3.90%  '
$ws.Cells.Item(31, 4).Value = '23.32'
$ws.Cells.Item(31, 5).Value = '  +54.88%  '
$ws.Cells.Item(32, 2).Value = 'Kaspa'
$ws.Cells.Item(32, 3).Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Cells.Item(32, 4).Value = '0.0978'
$ws.Cells.Item(32, 5).Value = '  +22.26%  '
$ws.Cells.Item(33, 2).Value = 'Filecoin'
$ws.Cells.Item(33, 3).Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Cells.Item(33, 4).Value = '4.50'
$ws.Cells.Item(33, 5).Value = '  -0.46%  '
$ws.Cells.Item(34, 2).Value = 'Hedera'
$ws.Cells.Item(34, 3).Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Cells.Item(34, 4).Value = '0.0602'
$ws.Cells.Item(34, 5).Value = '  -1.44%  '
$ws.Cells.Item(35, 2).Value = 'WEMIXToken'
$ws.Cells.Item(35, 3).Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Cells.Item(35, 4).Value = '1.91'
$ws.Cells.Item(35, 5).Value = '  -0.39%  '
$ws.Cells.Item(36, 5).Value = '  +0.00%  '
$ws.Cells.Item(37, 2).Value = 'LidoDAOToken'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Cells.Item(37, 4).Value = '2.30'
$ws.Cells.Item(37, 5).Value = '  +16.49%  '
$ws.Cells.Item(38, 4).Value = '4.11'
$ws.Cells.Item(38, 5).Value = '  -5.36%  '
$ws.Cells.Item(39, 4).Value = '0.898'
$ws.Cells.Item(39, 5).Value = '  +5.66%  '
$ws.Cells.Item(40, 4).Value = '1.35'
$ws.Cells.Item(40, 5).Value = '  -8.44%  '
$ws.Cells.Item(41, 4).Value = '1.18'
$ws.Cells.Item(41, 5).Value = '  +8.21%  '
$ws.Cells.Item(42, 2).Value = 'VeChain'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Cells.Item(42, 4).Value = '0.0221'
$ws.Cells.Item(42, 5).Value = '  -2.99%  '
$ws.Cells.Item(43, 2).Value = 'Aave'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Cells.Item(43, 4).Value = '98.60'
$ws.Cells.Item(43, 5).Value = '  -0.42%  '
$ws.Cells.Item(44, 4).Value = '2.81'
$ws.Cells.Item(44, 5).Value = '  +16.42%  '
$ws.Cells.Item(45, 4).Value = '16.27'
$ws.Cells.Item(45, 5).Value = '  -2.12%  '
$ws.Cells.Item(46, 4).Value = '1.366.37'
$ws.Cells.Item(46, 5).Value = '  +3.54%  '
$ws.Cells.Item(47, 4).Value = '0.0836'
$ws.Cells.Item(47, 5).Value = '  +3.31%  '
$ws.Cells.Item(48, 4).Value = '2.302.28'
$ws.Cells.Item(48, 5).Value = '  +10.97%  '
$ws.Cells.Item(49, 4).Value = '2.81'
$ws.Cells.Item(49, 5).Value = '  +2.07%  '
$ws.Cells.Item(50, 4).Value = '2.27'
$ws.Cells.Item(50, 5).Value = '  -2.74%  '
$ws.Cells.Item(51, 2).Value = 'FraxShare'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Cells.Item(51, 4).Value = '6.63'
$ws.Cells.Item(51, 5).Value = '  +4.57%  '

$priceVolRange.Style = "Normal"

Write-Host "Updated crypto table values."
